$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1314.8
$ws.Range("I6").Value = 574.75
$ws.Range("J6").Value = 1808.1666
$ws.Range("K6").Value = 1724.25
$ws.Range("L6").Value = 5424.4998
$ws.Range("M6").Value = -1612.25
$ws.Range("N6").Value = -5648.4998

$ws.Range("H12").Value = 457.6
$ws.Range("I12").Value = 529.3333
$ws.Range("J12").Value = 350
$ws.Range("K12").Value = 529.3333
$ws.Range("L12").Value = 350
$ws.Range("M12").Value = -359.3333
$ws.Range("N12").Value = -690

$ws.Range("H32").Value = 2508
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2508
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2508
$ws.Range("N32").Value = -3160

$ws.Range("H88").Value = 2299.375
$ws.Range("I88").Value = 1250
$ws.Range("J88").Value = 2649.1667
$ws.Range("K88").Value = 1250
$ws.Range("L88").Value = 2649.1667
$ws.Range("M88").Value = -844
$ws.Range("N88").Value = -3461.1667

$ws.Range("H91").Value = 2299.375
$ws.Range("I91").Value = 1250
$ws.Range("J91").Value = 2649.1667
$ws.Range("K91").Value = 1250
$ws.Range("L91").Value = 2649.1667
$ws.Range("M91").Value = 154
$ws.Range("N91").Value = -5457.1667

$ws.Range("H92").Value = 25000390
$ws.Range("I92").Value = 31250332
$ws.Range("J92").Value = 624
$ws.Range("K92").Value = 31250332
$ws.Range("L92").Value = 624
$ws.Range("M92").Value = -31249084
$ws.Range("N92").Value = -3120

$ws.Range("H98").Value = 2913.5
$ws.Range("I98").Value = 3049.9
$ws.Range("J98").Value = 2572.5
$ws.Range("K98").Value = 3049.9
$ws.Range("L98").Value = 2572.5
$ws.Range("M98").Value = -1551.9
$ws.Range("N98").Value = -5568.5

$ws.Range("H100").Value = 1629.5
$ws.Range("I100").Value = 966.1667
$ws.Range("J100").Value = 2624.5
$ws.Range("K100").Value = 966.1667
$ws.Range("L100").Value = 2624.5
$ws.Range("M100").Value = -425.1667
$ws.Range("N100").Value = -3706.5

$ws.Range("H116").Value = 15061.625
$ws.Range("I116").Value = 19083
$ws.Range("J116").Value = 2997.5
$ws.Range("K116").Value = 19083
$ws.Range("L116").Value = 2997.5
$ws.Range("M116").Value = -15641
$ws.Range("N116").Value = -9881.5

$ws.Range("H122").Value = 2913.5
$ws.Range("I122").Value = 3049.9
$ws.Range("J122").Value = 2572.5
$ws.Range("K122").Value = 9149.700000000001
$ws.Range("L122").Value = 7717.5
$ws.Range("M122").Value = -6699.700000000001
$ws.Range("N122").Value = -12617.5

$ws.Range("H129").Value = 1219.5518
$ws.Range("I129").Value = 889.3333
$ws.Range("J129").Value = 1257.6538
$ws.Range("K129").Value = 2667.9999
$ws.Range("L129").Value = 3772.9614
$ws.Range("M129").Value = 2332.0001
$ws.Range("N129").Value = -13772.9614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3537.4055
$ws.Range("I32").Value = 3046.9143
$ws.Range("J32").Value = 12121
$ws.Range("K32").Value = 3046.9143
$ws.Range("L32").Value = 12121
$ws.Range("M32").Value = -2759.9143
$ws.Range("N32").Value = -12695

$ws.Range("H45").Value = 1572.25
$ws.Range("I45").Value = 943.5
$ws.Range("J45").Value = 1886.625
$ws.Range("K45").Value = 943.5
$ws.Range("L45").Value = 1886.625
$ws.Range("M45").Value = -566.5
$ws.Range("N45").Value = -2640.625

$ws.Range("H132").Value = 1188.8636
$ws.Range("I132").Value = 892.3946999999999
$ws.Range("J132").Value = 3066.5
$ws.Range("K132").Value = 2677.1841
$ws.Range("L132").Value = 9199.5
$ws.Range("M132").Value = -147.1840999999999
$ws.Range("N132").Value = -14259.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1883.2222
$ws.Range("I94").Value = 328.5
$ws.Range("J94").Value = 4992.6665
$ws.Range("K94").Value = 328.5
$ws.Range("L94").Value = 4992.6665
$ws.Range("M94").Value = 122.5
$ws.Range("N94").Value = -5894.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1970.1
$ws.Range("I31").Value = 1812.125
$ws.Range("J31").Value = 2602
$ws.Range("K31").Value = 1812.125
$ws.Range("L31").Value = 2602
$ws.Range("M31").Value = -1517.125
$ws.Range("N31").Value = -3192

$ws.Range("H34").Value = 1970.1
$ws.Range("I34").Value = 1812.125
$ws.Range("J34").Value = 2602
$ws.Range("K34").Value = 1812.125
$ws.Range("L34").Value = 2602
$ws.Range("M34").Value = -1610.125
$ws.Range("N34").Value = -3006

$ws.Range("H58").Value = 1554656.9
$ws.Range("I58").Value = 1892022.1
$ws.Range("J58").Value = 2777
$ws.Range("K58").Value = 1892022.1
$ws.Range("L58").Value = 2777
$ws.Range("M58").Value = -1891819.1
$ws.Range("N58").Value = -3183

$ws.Range("H132").Value = 1366.0358
$ws.Range("I132").Value = 875.65216
$ws.Range("J132").Value = 3621.8
$ws.Range("K132").Value = 2626.95648
$ws.Range("L132").Value = 10865.4
$ws.Range("M132").Value = -96.95647999999983
$ws.Range("N132").Value = -15925.4

$ws.Range("H134").Value = 1754.4865
$ws.Range("I134").Value = 1591
$ws.Range("J134").Value = 3103.25
$ws.Range("K134").Value = 4773
$ws.Range("L134").Value = 9309.75
$ws.Range("M134").Value = -2238
$ws.Range("N134").Value = -14379.75

$ws.Range("H136").Value = 1554656.9
$ws.Range("I136").Value = 1892022.1
$ws.Range("J136").Value = 2777
$ws.Range("K136").Value = 5676066.300000001
$ws.Range("L136").Value = 8331
$ws.Range("M136").Value = -5673516.300000001
$ws.Range("N136").Value = -13431

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 259.2
$ws.Range("I6").Value = 259.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 777.5999999999999
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -664.5999999999999

$ws.Range("H7").Value = 718
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 739.8
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 2219.4
$ws.Range("M7").Value = -1388
$ws.Range("N7").Value = -2443.4

$ws.Range("H11").Value = 644.1429000000001
$ws.Range("I11").Value = 585
$ws.Range("J11").Value = 999
$ws.Range("K11").Value = 1755
$ws.Range("L11").Value = 2997
$ws.Range("M11").Value = -1615
$ws.Range("N11").Value = -3277

$ws.Range("H68").Value = 759.25
$ws.Range("I68").Value = 400
$ws.Range("J68").Value = 879
$ws.Range("K68").Value = 1200
$ws.Range("L68").Value = 2637
$ws.Range("M68").Value = -389
$ws.Range("N68").Value = -4259

$ws.Range("H71").Value = 759.25
$ws.Range("I71").Value = 400
$ws.Range("J71").Value = 879
$ws.Range("K71").Value = 3600
$ws.Range("L71").Value = 7911
$ws.Range("M71").Value = 456
$ws.Range("N71").Value = -16023

$ws.Range("H131").Value = 15991.674
$ws.Range("I131").Value = 707.5
$ws.Range("J131").Value = 17447.309
$ws.Range("K131").Value = 2122.5
$ws.Range("L131").Value = 52341.927
$ws.Range("M131").Value = 2917.5
$ws.Range("N131").Value = -62421.927

$ws.Range("H132").Value = 1780
$ws.Range("I132").Value = 1110
$ws.Range("J132").Value = 2896.6667
$ws.Range("K132").Value = 9990
$ws.Range("L132").Value = 26070.0003
$ws.Range("M132").Value = -7460
$ws.Range("N132").Value = -31130.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 15000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 15000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 15000
$ws.Range("N131").Value = -25080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2620.3076
$ws.Range("I7").Value = 2895.625
$ws.Range("J7").Value = 2179.8
$ws.Range("K7").Value = 2895.625
$ws.Range("L7").Value = 2179.8
$ws.Range("M7").Value = -2783.625
$ws.Range("N7").Value = -2403.8

$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 15000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -15452

$ws.Range("H82").Value = 1903.0667
$ws.Range("I82").Value = 1780
$ws.Range("J82").Value = 1964.6
$ws.Range("K82").Value = 1780
$ws.Range("L82").Value = 1964.6
$ws.Range("M82").Value = -1419
$ws.Range("N82").Value = -2686.6

$ws.Range("H85").Value = 1903.0667
$ws.Range("I85").Value = 1780
$ws.Range("J85").Value = 1964.6
$ws.Range("K85").Value = 1780
$ws.Range("L85").Value = 1964.6
$ws.Range("M85").Value = -532
$ws.Range("N85").Value = -4460.6

$ws.Range("H126").Value = 2620.3076
$ws.Range("I126").Value = 2895.625
$ws.Range("J126").Value = 2179.8
$ws.Range("K126").Value = 8686.875
$ws.Range("L126").Value = 6539.400000000001
$ws.Range("M126").Value = -6216.875
$ws.Range("N126").Value = -11479.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 14380
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 14380
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 14380
$ws.Range("N21").Value = -14850

$ws.Range("H30").Value = 15000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 15000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15214

$ws.Range("H35").Value = 14380
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 14380
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 14380
$ws.Range("N35").Value = -14960

$ws.Range("H81").Value = 1159.6
$ws.Range("I81").Value = 1159.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2319.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1258.2

$ws.Range("H84").Value = 1159.6
$ws.Range("I84").Value = 1159.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 11596
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -6292

$ws.Range("H107").Value = 626
$ws.Range("I107").Value = 544.8570999999999
$ws.Range("J107").Value = 788.2857
$ws.Range("K107").Value = 1634.5713
$ws.Range("L107").Value = 2364.8571
$ws.Range("M107").Value = 285.4287000000002
$ws.Range("N107").Value = -6204.8571

$ws.Range("H131").Value = 30000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 30000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080

$ws.Range("H136").Value = 25253988
$ws.Range("I136").Value = 34723310
$ws.Range("J136").Value = 2460
$ws.Range("K136").Value = 104169930
$ws.Range("L136").Value = 7380
$ws.Range("M136").Value = -104167380
$ws.Range("N136").Value = -12480
